# Generate Report for Archive
# - Flip the localization status text from "Ready for handoff" to
#   "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
#   de-de!C2 all share that status string).
# - Narrow the status column(s) that held that text: Overview cols E & F,
#   and column C on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column widths for the status columns ---
# (ColumnWidth 12.5 is the character-width input that lands on the same
# pixel-quantized column width the target file stores.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
